$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1301171.5
$ws.Range("I15").Value = 1301171.5
$ws.Range("K15").Value = 3903514.5
$ws.Range("M15").Value = -3903345.5

$ws.Range("H76").Value = 2749.6
$ws.Range("I76").Value = 2550.6428
$ws.Range("J76").Value = 3213.8333
$ws.Range("K76").Value = 2550.6428
$ws.Range("L76").Value = 3213.8333
$ws.Range("M76").Value = -2235.6428
$ws.Range("N76").Value = -3843.8333

$ws.Range("H79").Value = 2749.6
$ws.Range("I79").Value = 2550.6428
$ws.Range("J79").Value = 3213.8333
$ws.Range("K79").Value = 2550.6428
$ws.Range("L79").Value = 3213.8333
$ws.Range("M79").Value = -1458.6428
$ws.Range("N79").Value = -5397.8333

$ws.Range("H113").Value = 3260.0667
$ws.Range("I113").Value = 2743.889
$ws.Range("J113").Value = 4034.3333
$ws.Range("K113").Value = 2743.889
$ws.Range("L113").Value = 4034.3333
$ws.Range("M113").Value = 510.1109999999999
$ws.Range("N113").Value = -10542.3333

$ws.Range("H137").Value = 1555.06
$ws.Range("I137").Value = 2180
$ws.Range("J137").Value = 1064.0358
$ws.Range("K137").Value = 6540
$ws.Range("L137").Value = 3192.1074
$ws.Range("M137").Value = -3990
$ws.Range("N137").Value = -8292.107400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1012.6
$ws.Range("I2").Value = 622.13336
$ws.Range("J2").Value = 2184
$ws.Range("K2").Value = 622.13336
$ws.Range("L2").Value = 2184
$ws.Range("M2").Value = -509.13336
$ws.Range("N2").Value = -2410

$ws.Range("H45").Value = 1797.5385
$ws.Range("I45").Value = 1478
$ws.Range("J45").Value = 1939.5555
$ws.Range("K45").Value = 1478
$ws.Range("L45").Value = 1939.5555
$ws.Range("M45").Value = -1101
$ws.Range("N45").Value = -2693.5555

$ws.Range("H61").Value = 3352.6155
$ws.Range("I61").Value = 4545.909
$ws.Range("J61").Value = 1808.3529
$ws.Range("K61").Value = 4545.909
$ws.Range("L61").Value = 1808.3529
$ws.Range("M61").Value = -4333.909
$ws.Range("N61").Value = -2232.3529

$ws.Range("H63").Value = 4512.643
$ws.Range("J63").Value = 6262.4287
$ws.Range("L63").Value = 6262.4287
$ws.Range("N63").Value = -7634.4287

$ws.Range("H66").Value = 4512.643
$ws.Range("J66").Value = 6262.4287
$ws.Range("L66").Value = 31312.1435
$ws.Range("N66").Value = -38176.14350000001

$ws.Range("H116").Value = 1012.6
$ws.Range("I116").Value = 622.13336
$ws.Range("J116").Value = 2184
$ws.Range("K116").Value = 622.13336
$ws.Range("L116").Value = 2184
$ws.Range("M116").Value = 1671.86664
$ws.Range("N116").Value = -6772

$ws.Range("H122").Value = 1129
$ws.Range("I122").Value = 925
$ws.Range("J122").Value = 1265
$ws.Range("K122").Value = 2775
$ws.Range("L122").Value = 3795
$ws.Range("M122").Value = -325
$ws.Range("N122").Value = -8695

$ws.Range("H132").Value = 613972.1
$ws.Range("I132").Value = 1163713.4
$ws.Range("J132").Value = 4799.4053
$ws.Range("K132").Value = 3491140.2
$ws.Range("L132").Value = 14398.2159
$ws.Range("M132").Value = -3488610.2
$ws.Range("N132").Value = -19458.2159

$ws.Range("H136").Value = 3352.6155
$ws.Range("I136").Value = 4545.909
$ws.Range("J136").Value = 1808.3529
$ws.Range("K136").Value = 13637.727
$ws.Range("L136").Value = 5425.0587
$ws.Range("M136").Value = -11087.727
$ws.Range("N136").Value = -10525.0587

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1012.6
$ws.Range("I3").Value = 622.13336
$ws.Range("J3").Value = 2184
$ws.Range("K3").Value = 622.13336
$ws.Range("L3").Value = 2184
$ws.Range("M3").Value = -508.13336
$ws.Range("N3").Value = -2412

$ws.Range("H134").Value = 7377.1787
$ws.Range("I134").Value = 2460.0833
$ws.Range("J134").Value = 11065
$ws.Range("K134").Value = 7380.249899999999
$ws.Range("L134").Value = 33195
$ws.Range("M134").Value = -4845.249899999999
$ws.Range("N134").Value = -38265

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4071.4866
$ws.Range("I58").Value = 2465.7058
$ws.Range("J58").Value = 5436.4
$ws.Range("K58").Value = 2465.7058
$ws.Range("L58").Value = 5436.4
$ws.Range("M58").Value = -2262.7058
$ws.Range("N58").Value = -5842.4

$ws.Range("H122").Value = 200003300
$ws.Range("I122").Value = 500000400
$ws.Range("J122").Value = 5233.3335
$ws.Range("K122").Value = 1500001200
$ws.Range("L122").Value = 15700.0005
$ws.Range("M122").Value = -1499998750
$ws.Range("N122").Value = -20600.0005

$ws.Range("H136").Value = 4071.4866
$ws.Range("I136").Value = 2465.7058
$ws.Range("J136").Value = 5436.4
$ws.Range("K136").Value = 7397.117400000001
$ws.Range("L136").Value = 16309.2
$ws.Range("M136").Value = -4847.117400000001
$ws.Range("N136").Value = -21409.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 1297
$ws.Range("I103").Value = 292.5
$ws.Range("J103").Value = 1966.6666
$ws.Range("K103").Value = 877.5
$ws.Range("L103").Value = 5899.9998
$ws.Range("M103").Value = 1.5
$ws.Range("N103").Value = -7657.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 3000
$ws.Range("I43").Value = 3000
$ws.Range("K43").Value = 3000
$ws.Range("M43").Value = -2849

$ws.Range("H46").Value = 4516.6665
$ws.Range("I46").Value = 6425
$ws.Range("J46").Value = 700
$ws.Range("K46").Value = 6425
$ws.Range("L46").Value = 700
$ws.Range("M46").Value = -6269
$ws.Range("N46").Value = -1012

$ws.Range("H57").Value = 13995
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 13995
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 13995
$ws.Range("N57").Value = -15635
$ws.Range("M57").ClearContents()

$ws.Range("H102").Value = 2749245.5
$ws.Range("I102").Value = 4202912.5
$ws.Range("J102").Value = 3430.4443
$ws.Range("K102").Value = 4202912.5
$ws.Range("L102").Value = 3430.4443
$ws.Range("M102").Value = -4201290.5
$ws.Range("N102").Value = -6674.4443

$ws.Range("H122").Value = 2979.3845
$ws.Range("I122").Value = 2631
$ws.Range("J122").Value = 3637.4443
$ws.Range("K122").Value = 7893
$ws.Range("L122").Value = 10912.3329
$ws.Range("M122").Value = -5443
$ws.Range("N122").Value = -15812.3329

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1500
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 500
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = -205
$ws.Range("N22").Value = -2590

$ws.Range("H27").Value = 1500
$ws.Range("I27").Value = 500
$ws.Range("J27").Value = 2000
$ws.Range("K27").Value = 500
$ws.Range("L27").Value = 2000
$ws.Range("M27").Value = -393
$ws.Range("N27").Value = -2214

$ws.Range("H122").Value = 18172.715
$ws.Range("I122").Value = 27426
$ws.Range("J122").Value = 5835
$ws.Range("K122").Value = 82278
$ws.Range("L122").Value = 17505
$ws.Range("M122").Value = -79828
$ws.Range("N122").Value = -22405

$ws.Range("H132").Value = 82161.69500000001
$ws.Range("I132").Value = 171650.67
$ws.Range("J132").Value = 5456.857
$ws.Range("K132").Value = 514952.01
$ws.Range("L132").Value = 16370.571
$ws.Range("M132").Value = -512422.01
$ws.Range("N132").Value = -21430.571

$ws.Range("H136").Value = 2180.8
$ws.Range("I136").Value = 2228.5
$ws.Range("K136").Value = 6685.5
$ws.Range("M136").Value = -4135.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3145.0833
$ws.Range("I122").Value = 2104
$ws.Range("J122").Value = 4186.1665
$ws.Range("K122").Value = 6312
$ws.Range("L122").Value = 12558.4995
$ws.Range("M122").Value = -3862
$ws.Range("N122").Value = -17458.4995
